# Issue #15: Fix broken Drools tests.
# The fully-qualified DroolsActivity class name moved from
# com.centurylink.mdw.workflow.activity.rules.DroolsActivity
# to com.centurylink.mdw.drools.DroolsActivity -- update the
# "Import" cell on both decision-table sheets that reference it.

$wb = $excel.ActiveWorkbook

$newImport = "java.util.Map, java.util.Date, com.centurylink.mdw.drools.DroolsActivity"

$wsNorthSouth = $wb.Worksheets.Item("NorthSouth")
$wsNorthSouth.Range("C4").Value = $newImport

$wsEastWest = $wb.Worksheets.Item("EastWest")
$wsEastWest.Range("C4").Value = $newImport

# Restore the active-cell selections left behind on each sheet.
$wsEastWest.Activate() | Out-Null
$wsEastWest.Range("D25").Select() | Out-Null

$wsNorthSouth.Activate() | Out-Null
$wsNorthSouth.Range("D22").Select() | Out-Null
